# Applies the "introducing a scen file for GAMS commands example; creating
# commodity-only attributes case in SubRES" edit to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New commodity-only attributes block (rows 36-40, columns C:F) -------
# Cells are written in the exact order the new shared-string values first
# appear so the shared-string table is rebuilt with the same indices as the
# target file (TimeSlice, ELCHYD, S, Attribute, COM_BNDNET,
# "~FI_T: COM_FR~LO", COM_TAXNET).

$ws.Range("E37").Value = "TimeSlice"
$ws.Range("D38").Value = "ELCHYD"
$ws.Range("E38").Value = "S"
$ws.Range("C37").Value = "Attribute"
$ws.Range("C39").Value = "COM_BNDNET"
$ws.Range("E36").Value = "~FI_T: COM_FR~LO"
$ws.Range("C40").Value = "COM_TAXNET"

# Remaining cells reuse already-existing shared strings / plain numbers.
$ws.Range("D37").Value = "CommName"
$ws.Range("F37").Value = 2020
$ws.Range("F38").Value = 0.13
$ws.Range("D39").Value = "ELCSOL"
$ws.Range("F39").Value = 0
$ws.Range("D40").Value = "ELCWIN"
$ws.Range("F40").Value = 0

# --- View / selection -----------------------------------------------------
# Scroll the sheet so row 23 is at the top and leave the active cell at D41,
# matching the saved workbook view (topLeftCell="A23", selection D41).
[void]$ws.Activate()
[void]$ws.Range("A23").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D41").Select()
